$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: turn on turbulence (seed 100, Mann model) instead of "noturb"
$ws.Range("D4").Formula = '="dlc01_steady_wsp" & E4 & "_s100"'
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = 1
$ws.Range("J4").Value = "turb_s100_10ms"
$ws.Range("K4").Formula = "=E4*B4/512"

# Leave selection where the author last left it
$ws.Range("Q16").Select()
